# The source data was re-pulled/re-sorted upstream so that, within each
# 12-month year block, October/November/December now appear first,
# followed by January..September (values follow their original month
# label - only the row order within each year block changed).
#
# Rows 2-13  -> 2014 (Jan..Dec)
# Rows 14-25 -> 2015 (Jan..Dec)
# Rows 26-37 -> 2016 (Jan..Dec)
# Rows 38-49 -> 2017 (Jan..Dec)
#
# For each block, the new row order (0-based offsets into the original
# 12-row block) is: Oct, Nov, Dec, Jan, Feb, Mar, Apr, May, Jun, Jul, Aug, Sep

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$order = @(9,10,11,0,1,2,3,4,5,6,7,8)
$blockStarts = @(2,14,26,38)

foreach ($blockStart in $blockStarts) {
    $blockEnd = $blockStart + 11

    # Snapshot the 12 rows of this year block (columns A..F) before
    # overwriting anything, since rows get shuffled in place.
    $allRows = @()
    for ($r = $blockStart; $r -le $blockEnd; $r++) {
        $row = @(
            $ws.Cells.Item($r,1).Value2,
            $ws.Cells.Item($r,2).Value2,
            $ws.Cells.Item($r,3).Value2,
            $ws.Cells.Item($r,4).Value2,
            $ws.Cells.Item($r,5).Value2,
            $ws.Cells.Item($r,6).Value2
        )
        $allRows += ,$row
    }

    # Write the rows back in the new order.
    for ($i = 0; $i -lt 12; $i++) {
        $srcRow = $allRows[$order[$i]]
        $destR = $blockStart + $i
        for ($c = 1; $c -le 6; $c++) {
            $ws.Cells.Item($destR, $c).Value2 = $srcRow[$c - 1]
        }
    }
}
